$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 366; existing rows 366-413 shift down to 367-414.
$ws.Rows.Item(366).Insert()

# Populate the newly inserted row 366 with the new record.
$ws.Cells.Item(366, 1).Value = 10
$ws.Cells.Item(366, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(366, 3).Value = "La Araucanía"
$ws.Cells.Item(366, 4).Value = 45127
$ws.Cells.Item(366, 5).Value = 9
$ws.Cells.Item(366, 6).Value = 100112052
$ws.Cells.Item(366, 7).Value = "Albahaca"
$ws.Cells.Item(366, 8).Value = "Sin especificar"
$ws.Cells.Item(366, 9).Value = "Primera"
$ws.Cells.Item(366, 10).Value = 200
$ws.Cells.Item(366, 11).Value = 5000
$ws.Cells.Item(366, 12).Value = 6000
$ws.Cells.Item(366, 13).Value = 5250
$ws.Cells.Item(366, 14).Value = "$/paquete"
$ws.Cells.Item(366, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(366, 16).Value = 5250
$ws.Cells.Item(366, 17).Value = 1
$ws.Cells.Item(366, 18).Value = "Hortaliza"
